# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 08:08"

# 2) Pure numeric refreshes (no row shift) -----------------------------
# Row 27: Israel
$ws.Range("B27").Value = 273826
$ws.Range("C27").Value = 1517
$ws.Range("D27").Value = 208763
$ws.Range("E27").Value = 63306

# Row 59: Uzbekistan
$ws.Range("B59").Value = 59197
$ws.Range("C59").Value = 251
$ws.Range("D59").Value = 55751
$ws.Range("E59").Value = 2959
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 487

# Row 66: Kirguistan
$ws.Range("B66").Value = 47799
$ws.Range("C66").Value = 164
$ws.Range("D66").Value = 43644
$ws.Range("E66").Value = 3089

# Row 155: Belice
$ws.Range("B155").Value = 2204
$ws.Range("C155").Value = 8
$ws.Range("E155").Value = 796

# 3) Sri Lanka row: a new data row is inserted right after Somalia (row 139),
#    pushing Estonia/Tailandia/Gambia down by one row, and the old Sri Lanka
#    row (previously right before Malta) is removed.
$ws.Rows.Item(140).Insert()

$ws.Range("A140").Value = "Sri Lanka"
$ws.Range("B140").Value = 3733
$ws.Range("C140").Value = 220
$ws.Range("D140").Value = 3259
$ws.Range("E140").Value = 461
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 13

# After the insert, the old Sri Lanka row (originally row 143) has shifted
# down to row 144; remove it (it now sits just before Malta, row 145).
$ws.Rows.Item(144).Delete()
